$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1037.7384
$ws.Range("I17").Value = 581.5714
$ws.Range("K17").Value = 1744.7142
$ws.Range("M17").Value = -1576.7142
$ws.Range("H20").Value = 3000
$ws.Range("I20").Value = 3000
$ws.Range("K20").Value = 3000
$ws.Range("M20").Value = -2770
$ws.Range("H28").Value = 17604.727
$ws.Range("I28").Value = 5099.1665
$ws.Range("J28").Value = 32611.4
$ws.Range("K28").Value = 5099.1665
$ws.Range("L28").Value = 32611.4
$ws.Range("M28").Value = -4614.1665
$ws.Range("N28").Value = -33581.4
$ws.Range("H35").Value = 3000
$ws.Range("I35").Value = 3000
$ws.Range("K35").Value = 3000
$ws.Range("M35").Value = -2621
$ws.Range("H64").Value = 2927.7222
$ws.Range("J64").Value = 2981.7273
$ws.Range("L64").Value = 2981.7273
$ws.Range("N64").Value = -3477.7273
$ws.Range("H67").Value = 2927.7222
$ws.Range("J67").Value = 2981.7273
$ws.Range("L67").Value = 2981.7273
$ws.Range("N67").Value = -4697.7273
$ws.Range("H76").Value = 4632591
$ws.Range("I76").Value = 5053463
$ws.Range("K76").Value = 5053463
$ws.Range("M76").Value = -5053148
$ws.Range("H79").Value = 4632591
$ws.Range("I79").Value = 5053463
$ws.Range("K79").Value = 5053463
$ws.Range("M79").Value = -5052371
$ws.Range("H106").Value = 37503708
$ws.Range("I106").Value = 66671150
$ws.Range("J106").Value = 2714.2856
$ws.Range("K106").Value = 66671150
$ws.Range("L106").Value = 2714.2856
$ws.Range("M106").Value = -66670519
$ws.Range("N106").Value = -3976.2856
$ws.Range("H128").Value = 45918
$ws.Range("J128").Value = 45918
$ws.Range("L128").Value = 45918
$ws.Range("N128").Value = -55878
$ws.Range("H129").Value = 1369.5
$ws.Range("J129").Value = 1948.2941
$ws.Range("L129").Value = 5844.8823
$ws.Range("N129").Value = -15844.8823

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 924.61536
$ws.Range("I2").Value = 892.7273
$ws.Range("K2").Value = 892.7273
$ws.Range("M2").Value = -779.7273
$ws.Range("H32").Value = 1454842.9
$ws.Range("I32").Value = 1773174
$ws.Range("J32").Value = 22352.857
$ws.Range("K32").Value = 1773174
$ws.Range("L32").Value = 22352.857
$ws.Range("M32").Value = -1772887
$ws.Range("N32").Value = -22926.857
$ws.Range("H116").Value = 924.61536
$ws.Range("I116").Value = 892.7273
$ws.Range("K116").Value = 892.7273
$ws.Range("M116").Value = 1401.2727
$ws.Range("H134").Value = 34552.25
$ws.Range("J134").Value = 34552.25
$ws.Range("L134").Value = 34552.25
$ws.Range("N134").Value = -44692.25

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 924.61536
$ws.Range("I3").Value = 892.7273
$ws.Range("K3").Value = 892.7273
$ws.Range("M3").Value = -778.7273
$ws.Range("H86").Value = 1495.1818
$ws.Range("I86").Value = 1345.9459
$ws.Range("J86").Value = 2284
$ws.Range("K86").Value = 1345.9459
$ws.Range("L86").Value = 2284
$ws.Range("M86").Value = -222.9458999999999
$ws.Range("N86").Value = -4530
$ws.Range("H89").Value = 1495.1818
$ws.Range("I89").Value = 1345.9459
$ws.Range("J89").Value = 2284
$ws.Range("K89").Value = 6729.729499999999
$ws.Range("L89").Value = 11420
$ws.Range("M89").Value = -1113.729499999999
$ws.Range("N89").Value = -22652
$ws.Range("H105").Value = 12501197
$ws.Range("I105").Value = 12501197
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 12501197
$ws.Range("L105").Value = 0
$ws.Range("M105").ClearContents()
$ws.Range("N105").Value = -12499450
$ws.Range("H107").Value = 144771.28
$ws.Range("I107").Value = 168399.83
$ws.Range("K107").Value = 168399.83
$ws.Range("M107").Value = -166479.83

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 4300
$ws.Range("I62").Value = 4500
$ws.Range("K62").Value = 4500
$ws.Range("M62").Value = -3876
$ws.Range("H65").Value = 4300
$ws.Range("I65").Value = 4500
$ws.Range("K65").Value = 22500
$ws.Range("M65").Value = -19380
$ws.Range("H94").Value = 1714.625
$ws.Range("I94").Value = 0
$ws.Range("J94").Value = 1714.625
$ws.Range("K94").Value = 0
$ws.Range("L94").ClearContents()
$ws.Range("M94").Value = 1714.625
$ws.Range("N94").Value = -2616.625

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 6445.625
$ws.Range("I33").Value = 7900.615
$ws.Range("J33").Value = 140.66667
$ws.Range("K33").Value = 47403.69
$ws.Range("L33").Value = 844.0000200000001
$ws.Range("M33").Value = -47120.69
$ws.Range("N33").Value = -1410.00002
$ws.Range("H122").Value = 2231.6064
$ws.Range("I122").Value = 395.62964
$ws.Range("J122").Value = 3689.5881
$ws.Range("K122").Value = 3560.66676
$ws.Range("L122").Value = 33206.2929
$ws.Range("M122").Value = -1110.66676
$ws.Range("N122").Value = -38106.2929
$ws.Range("H125").Value = 4742
$ws.Range("J125").Value = 4742
$ws.Range("L125").Value = 14226
$ws.Range("N125").Value = -24066
$ws.Range("H126").Value = 4827.706
$ws.Range("J126").Value = 4827.706
$ws.Range("L126").Value = 14483.118
$ws.Range("N126").Value = -24363.118
$ws.Range("H137").Value = 29562.586
$ws.Range("I137").Value = 7475.1665
$ws.Range("K137").Value = 22425.4995
$ws.Range("M137").Value = -17325.4995

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H141").Value = 68135.28999999999
$ws.Range("J141").Value = 68135.28999999999
$ws.Range("L141").Value = 68135.28999999999
$ws.Range("N141").Value = -78495.28999999999

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 525.7143
$ws.Range("I55").Value = 517
$ws.Range("J55").Value = 528.0909
$ws.Range("K55").Value = 517
$ws.Range("L55").Value = 528.0909
$ws.Range("M55").Value = -344
$ws.Range("N55").Value = -874.0909
$ws.Range("H133").Value = 64826
$ws.Range("J133").Value = 64826
$ws.Range("L133").Value = 64826
$ws.Range("N133").Value = -69886
$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").ClearContents()
$ws.Range("N134").Value = 0

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 10115.857
$ws.Range("I45").Value = 7569
$ws.Range("J45").Value = 10540.333
$ws.Range("K45").Value = 7569
$ws.Range("L45").Value = 10540.333
$ws.Range("M45").Value = -7078
$ws.Range("N45").Value = -11522.333
$ws.Range("H105").Value = 90000
$ws.Range("J105").Value = 90000
$ws.Range("L105").Value = 90000
$ws.Range("N105").Value = -96988
$ws.Range("H107").Value = 351.76923
$ws.Range("I107").Value = 265.8889
$ws.Range("J107").Value = 545
$ws.Range("K107").Value = 797.6667
$ws.Range("L107").Value = 1635
$ws.Range("M107").Value = 1122.3333
$ws.Range("N107").Value = -5475
$ws.Range("H124").Value = 30357.25
$ws.Range("J124").Value = 30357.25
$ws.Range("L124").Value = 30357.25
$ws.Range("N124").Value = -40177.25
$ws.Range("H135").Value = 88431.25
$ws.Range("J135").Value = 88431.25
$ws.Range("L135").Value = 88431.25
$ws.Range("N135").Value = -98571.25
$ws.Range("H136").Value = 3515.5386
$ws.Range("I136").Value = 2895.2222
$ws.Range("J136").Value = 4911.25
$ws.Range("K136").Value = 8685.6666
$ws.Range("L136").Value = 14733.75
$ws.Range("M136").Value = -6135.6666
$ws.Range("N136").Value = -19833.75
$ws.Range("H138").Value = 74419.625
$ws.Range("J138").Value = 74419.625
$ws.Range("L138").Value = 74419.625
$ws.Range("N138").Value = -84699.625
